$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37731.4
$ws.Range("J3").Value = 37731.4
$ws.Range("L3").Value = 37731.4
$ws.Range("N3").Value = -37959.4
$ws.Range("H4").Value = 950
$ws.Range("I4").Value = 1187.5
$ws.Range("J4").Value = 475
$ws.Range("K4").Value = 1187.5
$ws.Range("L4").Value = 475
$ws.Range("M4").Value = -1073.5
$ws.Range("N4").Value = -703
$ws.Range("H5").Value = 45.57143
$ws.Range("I5").Value = 49
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 49
$ws.Range("L5").Value = 25
$ws.Range("M5").Value = 66
$ws.Range("N5").Value = -255
$ws.Range("H86").Value = 2236337.8
$ws.Range("J86").Value = 3680.2
$ws.Range("L86").Value = 3680.2
$ws.Range("N86").Value = -5926.2
$ws.Range("H89").Value = 2236337.8
$ws.Range("J89").Value = 3680.2
$ws.Range("L89").Value = 18401
$ws.Range("N89").Value = -29633
$ws.Range("H92").Value = 794.4545000000001
$ws.Range("I92").Value = 704.4
$ws.Range("K92").Value = 704.4
$ws.Range("M92").Value = 543.6
$ws.Range("H98").Value = 76924500
$ws.Range("I98").Value = 83334700
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 83334700
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -83333202
$ws.Range("N98").Value = -4996
$ws.Range("H102").Value = 37731.4
$ws.Range("J102").Value = 37731.4
$ws.Range("L102").Value = 37731.4
$ws.Range("N102").Value = -44221.4
$ws.Range("H105").Value = 51223.332
$ws.Range("J105").Value = 51223.332
$ws.Range("L105").Value = 51223.332
$ws.Range("N105").Value = -58211.332
$ws.Range("H112").Value = 908.625
$ws.Range("J112").Value = 1392.25
$ws.Range("L112").Value = 4176.75
$ws.Range("N112").Value = -6392.75
$ws.Range("H122").Value = 76924500
$ws.Range("I122").Value = 83334700
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 250004100
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -250001650
$ws.Range("N122").Value = -10900
$ws.Range("H129").Value = 2309
$ws.Range("I129").Value = 1078.6666
$ws.Range("K129").Value = 3235.9998
$ws.Range("M129").Value = 1764.0002
$ws.Range("H137").Value = 3184.2222
$ws.Range("I137").Value = 2171.6296
$ws.Range("J137").Value = 6222
$ws.Range("K137").Value = 6514.888800000001
$ws.Range("L137").Value = 18666
$ws.Range("M137").Value = -3964.888800000001
$ws.Range("N137").Value = -23766

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1073.6923
$ws.Range("I2").Value = 1073.6923
$ws.Range("K2").Value = 1073.6923
$ws.Range("M2").Value = -960.6922999999999
$ws.Range("H32").Value = 6759111
$ws.Range("I32").Value = 7144660
$ws.Range("K32").Value = 7144660
$ws.Range("M32").Value = -7144373
$ws.Range("H97").Value = 1099.76
$ws.Range("I97").Value = 1267.05
$ws.Range("J97").Value = 430.6
$ws.Range("K97").Value = 1267.05
$ws.Range("L97").Value = 430.6
$ws.Range("M97").Value = -771.05
$ws.Range("N97").Value = -1422.6
$ws.Range("H109").Value = 43922.125
$ws.Range("J109").Value = 43922.125
$ws.Range("L109").Value = 43922.125
$ws.Range("N109").Value = -46696.125
$ws.Range("H116").Value = 1073.6923
$ws.Range("I116").Value = 1073.6923
$ws.Range("K116").Value = 1073.6923
$ws.Range("M116").Value = 1220.3077

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1073.6923
$ws.Range("I3").Value = 1073.6923
$ws.Range("K3").Value = 1073.6923
$ws.Range("M3").Value = -959.6922999999999
$ws.Range("H94").Value = 996.1667
$ws.Range("I94").Value = 1052
$ws.Range("J94").Value = 956.2857
$ws.Range("K94").Value = 1052
$ws.Range("L94").Value = 956.2857
$ws.Range("M94").Value = -601
$ws.Range("N94").Value = -1858.2857
$ws.Range("H107").Value = 1506.6451
$ws.Range("I107").Value = 1560.7084
$ws.Range("J107").Value = 1321.2858
$ws.Range("K107").Value = 1560.7084
$ws.Range("L107").Value = 1321.2858
$ws.Range("M107").Value = 359.2916
$ws.Range("N107").Value = -5161.2858
$ws.Range("H117").Value = 117000
$ws.Range("J117").Value = 117000
$ws.Range("L117").Value = 117000
$ws.Range("N117").Value = -126178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 5184.6665
$ws.Range("J88").Value = 5621.6
$ws.Range("L88").Value = 5621.6
$ws.Range("N88").Value = -6433.6
$ws.Range("H91").Value = 5184.6665
$ws.Range("J91").Value = 5621.6
$ws.Range("L91").Value = 5621.6
$ws.Range("N91").Value = -8429.6
$ws.Range("H107").Value = 942.7
$ws.Range("I107").Value = 418.46155
$ws.Range("J107").Value = 1916.2858
$ws.Range("K107").Value = 418.46155
$ws.Range("L107").Value = 1916.2858
$ws.Range("M107").Value = 1501.53845
$ws.Range("N107").Value = -5756.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1665.1578
$ws.Range("I5").Value = 1658.6
$ws.Range("J5").Value = 1672.4445
$ws.Range("K5").Value = 4975.799999999999
$ws.Range("L5").Value = 5017.333500000001
$ws.Range("M5").Value = -4863.799999999999
$ws.Range("N5").Value = -5241.333500000001
$ws.Range("H14").Value = 557.6667
$ws.Range("I14").Value = 557.6667
$ws.Range("K14").Value = 1673.0001
$ws.Range("M14").Value = -1500.0001
$ws.Range("H56").Value = 5400
$ws.Range("I56").Value = 5400
$ws.Range("K56").Value = 5400
$ws.Range("M56").Value = -4870
$ws.Range("H132").Value = 1829.6666
$ws.Range("I132").Value = 1914.6364
$ws.Range("J132").Value = 1736.2
$ws.Range("K132").Value = 17231.7276
$ws.Range("L132").Value = 15625.8
$ws.Range("M132").Value = -14701.7276
$ws.Range("N132").Value = -20685.8
$ws.Range("H135").Value = 1665.1578
$ws.Range("I135").Value = 1658.6
$ws.Range("J135").Value = 1672.4445
$ws.Range("K135").Value = 14927.4
$ws.Range("L135").Value = 15052.0005
$ws.Range("M135").Value = -12392.4
$ws.Range("N135").Value = -20122.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3813.5
$ws.Range("I80").Value = 3701.8
$ws.Range("J80").Value = 3999.6667
$ws.Range("K80").Value = 3701.8
$ws.Range("L80").Value = 3999.6667
$ws.Range("M80").Value = -2703.8
$ws.Range("N80").Value = -5995.6667
$ws.Range("H83").Value = 3813.5
$ws.Range("I83").Value = 3701.8
$ws.Range("J83").Value = 3999.6667
$ws.Range("K83").Value = 18509
$ws.Range("L83").Value = 19998.3335
$ws.Range("M83").Value = -13517
$ws.Range("N83").Value = -29982.3335
$ws.Range("H104").Value = 45080.332
$ws.Range("J104").Value = 45080.332
$ws.Range("L104").Value = 45080.332
$ws.Range("N104").Value = -52068.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 706.25
$ws.Range("I16").Value = 619.1
$ws.Range("K16").Value = 619.1
$ws.Range("M16").Value = -449.1
$ws.Range("H22").Value = 3480
$ws.Range("I22").Value = 3475
$ws.Range("K22").Value = 3475
$ws.Range("M22").Value = -3180
$ws.Range("H27").Value = 3480
$ws.Range("I27").Value = 3475
$ws.Range("K27").Value = 3475
$ws.Range("M27").Value = -3368
$ws.Range("H46").Value = 5575.533
$ws.Range("I46").Value = 2104.125
$ws.Range("J46").Value = 9542.857
$ws.Range("K46").Value = 2104.125
$ws.Range("L46").Value = 9542.857
$ws.Range("M46").Value = -1916.125
$ws.Range("N46").Value = -9918.857
$ws.Range("H93").Value = 45461324
$ws.Range("I93").Value = 55563404
$ws.Range("K93").Value = 55563404
$ws.Range("M93").Value = -55562156
$ws.Range("H118").Value = 112000
$ws.Range("J118").Value = 112000
$ws.Range("L118").Value = 112000
$ws.Range("N118").Value = -115314
$ws.Range("H122").Value = 4383.825
$ws.Range("I122").Value = 4343
$ws.Range("K122").Value = 13029
$ws.Range("M122").Value = -10579

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 29488.5
$ws.Range("J41").Value = 29488.5
$ws.Range("L41").Value = 29488.5
$ws.Range("N41").Value = -30268.5
$ws.Range("H126").Value = 3888.3872
$ws.Range("I126").Value = 4110.125
$ws.Range("K126").Value = 12330.375
$ws.Range("M126").Value = -9860.375
